$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 984.7931
$ws.Range("J129").Value = 1240
$ws.Range("L129").Value = 3720
$ws.Range("N129").Value = -13720
$ws.Range("H132").Value = 239711.6
$ws.Range("I132").Value = 1638.6774
$ws.Range("J132").Value = 910644.4
$ws.Range("K132").Value = 4916.0322
$ws.Range("L132").Value = 2731933.2
$ws.Range("M132").Value = -2386.0322
$ws.Range("N132").Value = -2736993.2
$ws.Range("H137").Value = 1380
$ws.Range("I137").Value = 1166.6666
$ws.Range("J137").Value = 1806.6666
$ws.Range("K137").Value = 3499.9998
$ws.Range("L137").Value = 5419.9998
$ws.Range("M137").Value = -949.9998000000001
$ws.Range("N137").Value = -10519.9998
$ws.Range("H138").Value = 3251.37
$ws.Range("I138").Value = 783.8570999999999
$ws.Range("J138").Value = 3907.2913
$ws.Range("K138").Value = 2351.5713
$ws.Range("L138").Value = 11721.8739
$ws.Range("M138").Value = 2788.4287
$ws.Range("N138").Value = -22001.8739

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 6000
$ws.Range("J46").Value = 6000
$ws.Range("L46").Value = 6000
$ws.Range("N46").Value = -6638
$ws.Range("H61").Value = 1453.475
$ws.Range("I61").Value = 1300.6333
$ws.Range("J61").Value = 1912
$ws.Range("K61").Value = 1300.6333
$ws.Range("L61").Value = 1912
$ws.Range("M61").Value = -1088.6333
$ws.Range("N61").Value = -2336
$ws.Range("H74").Value = 780.8982999999999
$ws.Range("I74").Value = 774.3555
$ws.Range("J74").Value = 801.9286
$ws.Range("K74").Value = 774.3555
$ws.Range("L74").Value = 801.9286
$ws.Range("M74").Value = 99.64449999999999
$ws.Range("N74").Value = -2549.9286
$ws.Range("H77").Value = 780.8982999999999
$ws.Range("I77").Value = 774.3555
$ws.Range("J77").Value = 801.9286
$ws.Range("K77").Value = 3871.7775
$ws.Range("L77").Value = 4009.643
$ws.Range("M77").Value = 496.2224999999999
$ws.Range("N77").Value = -12745.643
$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344
$ws.Range("H105").Value = 39400
$ws.Range("J105").Value = 39400
$ws.Range("L105").Value = 39400
$ws.Range("N105").Value = -46388
$ws.Range("H118").Value = 31533.334
$ws.Range("J118").Value = 31533.334
$ws.Range("L118").Value = 31533.334
$ws.Range("N118").Value = -34847.334
$ws.Range("H132").Value = 1494.0303
$ws.Range("I132").Value = 939
$ws.Range("J132").Value = 2604.0908
$ws.Range("K132").Value = 2817
$ws.Range("L132").Value = 7812.2724
$ws.Range("M132").Value = -287
$ws.Range("N132").Value = -12872.2724
$ws.Range("H136").Value = 1453.475
$ws.Range("I136").Value = 1300.6333
$ws.Range("J136").Value = 1912
$ws.Range("K136").Value = 3901.8999
$ws.Range("L136").Value = 5736
$ws.Range("M136").Value = -1351.8999
$ws.Range("N136").Value = -10836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 74716
$ws.Range("I134").Value = 3564.3333
$ws.Range("K134").Value = 10692.9999
$ws.Range("M134").Value = -8157.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H31").Value = 2710.739
$ws.Range("I31").Value = 2752.1365
$ws.Range("J31").Value = 1800
$ws.Range("K31").Value = 2752.1365
$ws.Range("L31").Value = 1800
$ws.Range("M31").Value = -2457.1365
$ws.Range("N31").Value = -2390
$ws.Range("H34").Value = 2710.739
$ws.Range("I34").Value = 2752.1365
$ws.Range("J34").Value = 1800
$ws.Range("K34").Value = 2752.1365
$ws.Range("L34").Value = 1800
$ws.Range("M34").Value = -2550.1365
$ws.Range("N34").Value = -2204
$ws.Range("H58").Value = 6921.1
$ws.Range("I58").Value = 2359.6
$ws.Range("J58").Value = 11482.6
$ws.Range("K58").Value = 2359.6
$ws.Range("L58").Value = 11482.6
$ws.Range("M58").Value = -2156.6
$ws.Range("N58").Value = -11888.6
$ws.Range("H62").Value = 9440.333000000001
$ws.Range("I62").Value = 11075.625
$ws.Range("J62").Value = 7571.4287
$ws.Range("K62").Value = 11075.625
$ws.Range("L62").Value = 7571.4287
$ws.Range("M62").Value = -10451.625
$ws.Range("N62").Value = -8819.4287
$ws.Range("H65").Value = 9440.333000000001
$ws.Range("I65").Value = 11075.625
$ws.Range("J65").Value = 7571.4287
$ws.Range("K65").Value = 55378.125
$ws.Range("L65").Value = 37857.14350000001
$ws.Range("M65").Value = -52258.125
$ws.Range("N65").Value = -44097.14350000001
$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180
$ws.Range("H132").Value = 1947.1794
$ws.Range("I132").Value = 1367.88
$ws.Range("J132").Value = 2981.6428
$ws.Range("K132").Value = 4103.64
$ws.Range("L132").Value = 8944.928400000001
$ws.Range("M132").Value = -1573.64
$ws.Range("N132").Value = -14004.9284
$ws.Range("H134").Value = 2166.875
$ws.Range("I134").Value = 1725.6666
$ws.Range("J134").Value = 3009.182
$ws.Range("K134").Value = 5176.9998
$ws.Range("L134").Value = 9027.545999999998
$ws.Range("M134").Value = -2641.9998
$ws.Range("N134").Value = -14097.546
$ws.Range("H136").Value = 6921.1
$ws.Range("I136").Value = 2359.6
$ws.Range("J136").Value = 11482.6
$ws.Range("K136").Value = 7078.799999999999
$ws.Range("L136").Value = 34447.8
$ws.Range("M136").Value = -4528.799999999999
$ws.Range("N136").Value = -39547.8
$ws.Range("H140").Value = 52476.363
$ws.Range("J140").Value = 52476.363
$ws.Range("L140").Value = 52476.363
$ws.Range("N140").Value = -62836.363

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 16992.309
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 16992.309
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 50976.927
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -60652.927

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8500
$ws.Range("J5").Value = 8500
$ws.Range("L5").Value = 8500
$ws.Range("N5").Value = -8724
$ws.Range("H70").Value = 4300.4736
$ws.Range("I70").Value = 3838.0715
$ws.Range("J70").Value = 5595.2
$ws.Range("K70").Value = 3838.0715
$ws.Range("L70").Value = 5595.2
$ws.Range("M70").Value = -3568.0715
$ws.Range("N70").Value = -6135.2
$ws.Range("H73").Value = 4300.4736
$ws.Range("I73").Value = 3838.0715
$ws.Range("J73").Value = 5595.2
$ws.Range("K73").Value = 3838.0715
$ws.Range("L73").Value = 5595.2
$ws.Range("M73").Value = -2902.0715
$ws.Range("N73").Value = -7467.2
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 2039.8684
$ws.Range("I122").Value = 1488
$ws.Range("J122").Value = 2985.9285
$ws.Range("K122").Value = 4464
$ws.Range("L122").Value = 8957.7855
$ws.Range("M122").Value = -2014
$ws.Range("N122").Value = -13857.7855
$ws.Range("H132").Value = 2981.8823
$ws.Range("I132").Value = 2426.7273
$ws.Range("K132").Value = 7280.1819
$ws.Range("M132").Value = -4750.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2699.7693
$ws.Range("I132").Value = 2141.8948
$ws.Range("K132").Value = 6425.6844
$ws.Range("M132").Value = -3895.6844
$ws.Range("H136").Value = 2873
$ws.Range("I136").Value = 1235.8096
$ws.Range("J136").Value = 14333.333
$ws.Range("K136").Value = 3707.4288
$ws.Range("L136").Value = 42999.999
$ws.Range("M136").Value = -1157.4288
$ws.Range("N136").Value = -48099.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H132").Value = 1616.907
$ws.Range("I132").Value = 1264.5151
$ws.Range("K132").Value = 3793.5453
$ws.Range("M132").Value = -1263.5453
$ws.Range("H136").Value = 1394.2368
$ws.Range("I136").Value = 1412.0883
$ws.Range("K136").Value = 4236.2649
$ws.Range("M136").Value = -1686.2649
